# Auto-generated edit script applying Adamantoise_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1113.4
$ws.Range("I33").Value = 767
$ws.Range("K33").Value = 767
$ws.Range("M33").Value = -538

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H76").Value = 4124.25
$ws.Range("I76").Value = 3999
$ws.Range("K76").Value = 3999
$ws.Range("M76").Value = -3684

$ws.Range("H79").Value = 4124.25
$ws.Range("I79").Value = 3999
$ws.Range("K79").Value = 3999
$ws.Range("M79").Value = -2907

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2944.5
$ws.Range("I2").Value = 2901.7693
$ws.Range("K2").Value = 2901.7693
$ws.Range("M2").Value = -2788.7693

$ws.Range("H45").Value = 5354.9287
$ws.Range("I45").Value = 5107.778
$ws.Range("J45").Value = 5799.8
$ws.Range("K45").Value = 5107.778
$ws.Range("L45").Value = 5799.8
$ws.Range("M45").Value = -4730.778
$ws.Range("N45").Value = -6553.8

$ws.Range("H116").Value = 2944.5
$ws.Range("I116").Value = 2901.7693
$ws.Range("K116").Value = 2901.7693
$ws.Range("M116").Value = -607.7692999999999

$ws.Range("H132").Value = 2540.7708
$ws.Range("I132").Value = 2277.5676
$ws.Range("K132").Value = 6832.702799999999
$ws.Range("M132").Value = -4302.702799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2944.5
$ws.Range("I3").Value = 2901.7693
$ws.Range("K3").Value = 2901.7693
$ws.Range("M3").Value = -2787.7693

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H94").Value = 1043.931
$ws.Range("I94").Value = 828.5294
$ws.Range("J94").Value = 1349.0834
$ws.Range("K94").Value = 828.5294
$ws.Range("L94").Value = 1349.0834
$ws.Range("M94").Value = -377.5294
$ws.Range("N94").Value = -2251.0834

$ws.Range("H99").Value = 4120.5386
$ws.Range("I99").Value = 3002.5
$ws.Range("K99").Value = 3002.5
$ws.Range("M99").Value = -1504.5

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1873.4
$ws.Range("I16").Value = 1190.5714
$ws.Range("J16").Value = 3466.6667
$ws.Range("K16").Value = 1190.5714
$ws.Range("L16").Value = 3466.6667
$ws.Range("M16").Value = -903.5714
$ws.Range("N16").Value = -4040.6667

$ws.Range("H23").Value = 19009
$ws.Range("I23").Value = 19009
$ws.Range("K23").Value = 19009
$ws.Range("M23").Value = -18769

$ws.Range("H27").Value = 19009
$ws.Range("I27").Value = 19009
$ws.Range("K27").Value = 19009
$ws.Range("M27").Value = -18817

$ws.Range("H96").Value = 39466.668
$ws.Range("J96").Value = 39466.668
$ws.Range("L96").Value = 39466.668
$ws.Range("N96").Value = -44958.668

$ws.Range("H113").Value = 1873.4
$ws.Range("I113").Value = 1190.5714
$ws.Range("J113").Value = 3466.6667
$ws.Range("K113").Value = 1190.5714
$ws.Range("L113").Value = 3466.6667
$ws.Range("M113").Value = 979.4286
$ws.Range("N113").Value = -7806.6667

$ws.Range("H132").Value = 3557
$ws.Range("I132").Value = 3415.0605
$ws.Range("K132").Value = 10245.1815
$ws.Range("M132").Value = -7715.181500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3674.875
$ws.Range("I3").Value = 4979.8
$ws.Range("K3").Value = 14939.4
$ws.Range("M3").Value = -14827.4

$ws.Range("H7").Value = 767.8378
$ws.Range("J7").Value = 894
$ws.Range("L7").Value = 2682
$ws.Range("N7").Value = -2906

$ws.Range("H9").Value = 4999.25
$ws.Range("I9").Value = 4998
$ws.Range("K9").Value = 14994
$ws.Range("M9").Value = -14770

$ws.Range("H18").Value = 6276.3335
$ws.Range("I18").Value = 4830
$ws.Range("K18").Value = 14490
$ws.Range("M18").Value = -14321

$ws.Range("H134").Value = 5791.9287
$ws.Range("I134").Value = 5791.9287
$ws.Range("K134").Value = 17375.7861
$ws.Range("M134").Value = -12305.7861

$ws.Range("H139").Value = 2282.9
$ws.Range("I139").Value = 2297.7896
$ws.Range("K139").Value = 6893.3688
$ws.Range("M139").Value = -1753.3688

$ws.Range("H140").Value = 1901.7778
$ws.Range("I140").Value = 1648.8
$ws.Range("K140").Value = 4946.4
$ws.Range("M140").Value = 233.6000000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -361
$ws.Range("N13").ClearContents()

$ws.Range("H132").Value = 3594.1936
$ws.Range("I132").Value = 3051.5833
$ws.Range("K132").Value = 9154.749899999999
$ws.Range("M132").Value = -6624.749899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 11010.5
$ws.Range("I34").Value = 8021
$ws.Range("K34").Value = 8021
$ws.Range("M34").Value = -7849

$ws.Range("H132").Value = 6552.815
$ws.Range("I132").Value = 6122.375
$ws.Range("K132").Value = 18367.125
$ws.Range("M132").Value = -15837.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 670.1667
$ws.Range("I107").Value = 506
$ws.Range("K107").Value = 1518
$ws.Range("M107").Value = 402

$ws.Range("H113").Value = 953.1
$ws.Range("I113").Value = 1323.6666
$ws.Range("J113").Value = 397.25
$ws.Range("K113").Value = 3970.9998
$ws.Range("L113").Value = 1191.75
$ws.Range("M113").Value = -1800.9998
$ws.Range("N113").Value = -5531.75

$ws.Range("H122").Value = 76928400
$ws.Range("I122").Value = 333335170
$ws.Range("J122").Value = 6375.1
$ws.Range("K122").Value = 1000005510
$ws.Range("L122").Value = 19125.3
$ws.Range("M122").Value = -1000003060
$ws.Range("N122").Value = -24025.3

$ws.Range("H132").Value = 1891.8549
$ws.Range("I132").Value = 1758.9344
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 5276.8032
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -2746.8032
$ws.Range("N132").Value = -35060
